# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Adds "Latest Target File" (F) / "Latest Handback File" (G) columns' data + hyperlinks
#    for the zh-cn and de-de detail sheets
#  - Stamps the "Latest Handback DateTime" (H) for each language

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Cornflower blue (FF6495ED), matching the workbook's existing HyperLink style,
# expressed as the VBA-style BGR decimal that Font.Color expects.
$hyperlinkColor = 15570276

function Apply-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Overview sheet: refresh the Status columns (B = zh-cn, C = de-de)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$mdName = "92fdf762-515f-46c8-b36b-7e203777c5d1.md"
$zhcnXlfName = "92fdf762-515f-46c8-b36b-7e203777c5d1.684d96098a4e7263fc3c9c5a64347c9855701a2f.zh-cn.xlf"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0dad6c8973e1cab511d21457722cc66f1265167e/e2e/92fdf762-515f-46c8-b36b-7e203777c5d1.md"
$zhcnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc1676b698a261a7ad03d137e13ef19c77aac4f8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/92fdf762-515f-46c8-b36b-7e203777c5d1.684d96098a4e7263fc3c9c5a64347c9855701a2f.zh-cn.xlf"

$zhcn.Range("F2").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $mdUrl, "", "", $mdName)
Apply-HyperlinkLook $zhcn.Range("F2")

$zhcn.Range("G2").Value = $zhcnXlfName
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zhcnXlfUrl, "", "", $zhcnXlfName)
Apply-HyperlinkLook $zhcn.Range("G2")

$zhcn.Range("F3").Value = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $mdUrl, "", "", $mdName)
Apply-HyperlinkLook $zhcn.Range("F3")

$zhcn.Range("G3").Value = $zhcnXlfName
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhcnXlfUrl, "", "", $zhcnXlfName)
Apply-HyperlinkLook $zhcn.Range("G3")

$zhcn.Range("H2").Value = "2016-03-19 02:49:31"
$zhcn.Range("H3").Value = "2016-03-19 02:49:31"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dedeXlfName = "92fdf762-515f-46c8-b36b-7e203777c5d1.684d96098a4e7263fc3c9c5a64347c9855701a2f.de-de.xlf"
$dedeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1616c4b1e3d8b237526f49487d0f4e8fbf2b9e5c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/92fdf762-515f-46c8-b36b-7e203777c5d1.684d96098a4e7263fc3c9c5a64347c9855701a2f.de-de.xlf"

$dede.Range("F2").Value = $mdName
$dede.Hyperlinks.Add($dede.Range("F2"), $mdUrl, "", "", $mdName)
Apply-HyperlinkLook $dede.Range("F2")

$dede.Range("G2").Value = $dedeXlfName
$dede.Hyperlinks.Add($dede.Range("G2"), $dedeXlfUrl, "", "", $dedeXlfName)
Apply-HyperlinkLook $dede.Range("G2")

$dede.Range("F3").Value = $mdName
$dede.Hyperlinks.Add($dede.Range("F3"), $mdUrl, "", "", $mdName)
Apply-HyperlinkLook $dede.Range("F3")

$dede.Range("G3").Value = $dedeXlfName
$dede.Hyperlinks.Add($dede.Range("G3"), $dedeXlfUrl, "", "", $dedeXlfName)
Apply-HyperlinkLook $dede.Range("G3")

$dede.Range("H2").Value = "2016-03-19 02:49:36"
$dede.Range("H3").Value = "2016-03-19 02:49:36"
